# Insert a new weekly record at row 52 (shifting existing rows 52..118 down
# to 53..119) on the single "Sheet1" worksheet, then populate the newly
# inserted row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52..118 down by one row, growing the used range to row 119.
$ws.Rows.Item(52).Insert()

# Populate the new row 52 with the new data record.
$ws.Cells.Item(52, 1).Value = 10
$ws.Cells.Item(52, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(52, 3).Value = "La Araucanía"
$ws.Cells.Item(52, 4).Value = 44413
$ws.Cells.Item(52, 5).Value = 9
$ws.Cells.Item(52, 6).Value = 100112039
$ws.Cells.Item(52, 7).Value = "Ciboulette"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 20
$ws.Cells.Item(52, 11).Value = 10000
$ws.Cells.Item(52, 12).Value = 10000
$ws.Cells.Item(52, 13).Value = 10000
$ws.Cells.Item(52, 14).Value = "$/docena de atados"
$ws.Cells.Item(52, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(52, 16).Value = 3333
$ws.Cells.Item(52, 17).Value = 3
$ws.Cells.Item(52, 18).Value = "Hortaliza"
